$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "1.00" / "0.527" / "0.000285" that Excel
# would otherwise auto-convert to a number on assignment. A leading apostrophe
# (the same trick used when typing numeric-looking text directly into a cell)
# keeps those values stored as plain text, matching the rest of the sheet.

$ws.Range("D2").Value = "68.116.55"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.782.29"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.90"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'170.59"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "3.782.50"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "'0.0000282"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "'36.76"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "4.417.90"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "3.788.98"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "'18.91"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "68.108.71"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "'10.64"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'468.47"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'0.0000151"
$ws.Range("E24").Value = "  -7.55%  "
$ws.Range("D25").Value = "'83.84"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "'12.17"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'10.57"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").Value = "3.932.18"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "'7.62"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "'30.57"
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'9.26"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "3.745.03"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.77"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").Value = "'5.86"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'0.315"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "'8.71"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.96"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "'404.46"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000285"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'45.70"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "'40.39"
$ws.Range("E50").Value = "  +7.72%  "
$ws.Range("D51").Value = "'140.85"
$ws.Range("E51").Value = "  -0.59%  "
